# Generate Report for Handback
# Update the timestamp strings recorded for the "37ea0299-ae78-4fad-9192-93c0c39d1158.md"
# file's handoff/handback generation times across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 37ea0299-...md
$wsOverview.Range("G3").Value = "2016-08-20 16:54:25"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for 37ea0299-...2484b996eabb1551a5ed41e872e80d74ed90610d.zh-cn.xlf
$wsZhCn.Range("H3").Value = "2016-08-20 16:54:21"
$wsZhCn.Range("K3").Value = "2016-08-20 16:54:39"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview value) and
# "Correspond Handback DateTime" for 37ea0299-...2484b996eabb1551a5ed41e872e80d74ed90610d.de-de.xlf
$wsDeDe.Range("H3").Value = "2016-08-20 16:54:25"
$wsDeDe.Range("K3").Value = "2016-08-20 16:54:45"
